$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Updated "Data" (GDP per Capita) values for existing years 1950-2010 (rows 2-62)
$updatedValues = @{
    2 = "861"
    3 = "862"
    4 = "874"
    5 = "872"
    6 = "872"
    7 = "810"
    8 = "878"
    9 = "845"
    10 = "813"
    11 = "838"
    12 = "869"
    13 = "899"
    14 = "877"
    15 = "947"
    16 = "939"
    17 = "969"
    18 = "963"
    19 = "921"
    20 = "987"
    21 = "980"
    22 = "1004"
    23 = "934"
    24 = "807"
    25 = "792"
    26 = "872"
    27 = "843"
    28 = "862"
    29 = "843"
    30 = "878"
    31 = "894"
    32 = "875"
    33 = "878"
    34 = "869"
    35 = "891"
    36 = "918"
    37 = "925"
    38 = "944"
    39 = "960"
    40 = "964"
    41 = "966"
    42 = "1006"
    43 = "1031.42994915762"
    44 = "1062.46039963364"
    45 = "1095.65638403607"
    46 = "1136.41846534845"
    47 = "1182.95140041245"
    48 = "1234.08406379444"
    49 = "1291.5125695079"
    50 = "1348.1499270396"
    51 = "1412.98132124769"
    52 = "1485.30922045797"
    53 = "1551.19647671542"
    54 = "1621.26697138793"
    55 = "1710.73135124973"
    56 = "1812.62466142731"
    57 = "1925.60616962229"
    58 = "2057.698392171"
    59 = "2194.42233788701"
    60 = "2321.05045875432"
    61 = "2448.88687843863"
    62 = "2599.2083634666"
}

foreach ($row in $updatedValues.Keys) {
    $cell = $ws.Cells.Item($row, 5)
    $cell.NumberFormat = "@"
    $cell.Value = $updatedValues[$row]
    $cell.ClearFormats()
}

# New rows for years 2011-2016 (rows 63-68)
$newRows = @{
    63 = @(2011, "2772")
    64 = @(2012, "2914")
    65 = @(2013, "3057")
    66 = @(2014, "3217")
    67 = @(2015, "3402")
    68 = @(2016, "3604")
}

foreach ($row in $newRows.Keys) {
    $year = $newRows[$row][0]
    $val = $newRows[$row][1]
    $ws.Cells.Item($row, 1).Value = 50
    $ws.Cells.Item($row, 2).Value = "Bangladesh"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $year
    $valCell = $ws.Cells.Item($row, 5)
    $valCell.NumberFormat = "@"
    $valCell.Value = $val
    $valCell.ClearFormats()
}
